# Add a new "free" male slot (ゆうた) at the top of the base_parameter
# roster and a second "ともみ" entry further down; every row that used to
# sit below those insertion points shifts down by one (two rows total),
# which is why rows 18-34 all end up with the data that used to be one
# (or two) rows above them, and the two previously-empty buffer rows
# (33, 34) get populated with the data that spilled into them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("base_parameter")

# Columns A..F for rows 18..34 in their final (post-edit) state.
$colA = @("〇", "",  "",  "",  "",  "",   "〇", "",  "",  "",  "〇", "",  "〇", "",  "",  "",  "")
$colB = @("ゆうた", "[FREE_W1]", "[FREE_W2]", "[FREE_W3]", "[FREE_W4]", "かえで", "シャンシャン", "たぬき", "ともね", "ともみ", "ともみ", "ななみ", "ひろみ", "まや", "みく", "れいか", "わかな")
$colC = @("男", "女", "女", "女", "女", "女", "女", "女", "女", "女", "女", "女", "女", "女", "女", "女", "女")
$colD = @(4, 1, 2, 3, 4, 3.5, 1, 2.5, 2.5, 3.5, 3.5, 2.5, 3.5, 3.5, 3, 3.5, 2)
$colE = @("", "", "", "", "", "〇", "", "", "", "", "", "", "", "", "", "", "")
$colF = @("〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇", "〇")

$startRow = 18
for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
    $ws.Cells.Item($r, 5).Value = $colE[$i]
    $ws.Cells.Item($r, 6).Value = $colF[$i]
}

# The author's last recorded selection on this sheet moved to E8.
$ws.Range("E8").Select() | Out-Null
